{"js": "// Applies the textual corrections described in the commit diff to the\n// synopsis document. Each change is performed with a body.search() +\n// insertText(..., Word.InsertLocation.replace) pair so the edit is\n// anchored on unambiguous surrounding context.\n\nasync function replaceOnce(context, searchText, newText, options) {\n  const body = context.document.body;\n  const searchOptions = Object.assign({ matchCase: true }, options || {});\n  const results = body.search(searchText, searchOptions);\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + searchText);\n  }\n\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 1. \"ressource\" -> \"ressources\"\nawait replaceOnce(\n  context,\n  \"un nombre limit\u00e9 de ressource qu\u2019il devra g\u00e9rer\",\n  \"un nombre limit\u00e9 de ressources qu\u2019il devra g\u00e9rer\"\n);\n\n// 2. add commas around \"qui servira pour bon nombre de choses\"\nawait replaceOnce(\n  context,\n  \"g\u00e9rer son budget qui servira pour bon nombre de choses \u00e0 savoir marchander\",\n  \"g\u00e9rer son budget, qui servira pour bon nombre de choses, \u00e0 savoir marchander\"\n);\n\n// 3. add comma after \"croissante\"\nawait replaceOnce(\n  context,\n  \"le logement d\u2019une population croissante ou encore pour la gestion\",\n  \"le logement d\u2019une population croissante, ou encore pour la gestion\"\n);\n\n// 4. remove duplicate \"l\u00e0\"\nawait replaceOnce(\n  context,\n  \"mais \u00e9galement l\u00e0 pour vous racheter\",\n  \"mais \u00e9galement pour vous racheter\"\n);\n\n// 5. \"tel que\" -> \"telles que\" + commas\nawait replaceOnce(\n  context,\n  \"au cours du jeu tel que l\u2019arbre de comp\u00e9tences octroyant des bonus non n\u00e9gligeables ou encore un syst\u00e8me\",\n  \"au cours du jeu, telles que l\u2019arbre de comp\u00e9tences octroyant des bonus non n\u00e9gligeables, ou encore un syst\u00e8me\"\n);\n\n// 6. \"car dans\" -> \"car, dans\"\nawait replaceOnce(\n  context,\n  \" car dans le cas contraire cela veut dire\",\n  \" car, dans le cas contraire cela veut dire\"\n);\n\n// 7. \"ne doivent pas\" -> \"de\"\nawait replaceOnce(\n  context,\n  \"Emp\u00eacher toutes les IA de la carte ne doivent pas faire la m\u00eame chose\",\n  \"Emp\u00eacher toutes les IA de la carte de faire la m\u00eame chose\"\n);\n\n// 8. \"etc\" -> \"etc.\"\nawait replaceOnce(\n  context,\n  \"ressources etc + arbres de comp\u00e9tences\",\n  \"ressources etc. + arbres de comp\u00e9tences\"\n);\n\n// 9. \"ces arbres\" -> \"les arbres\"\nawait replaceOnce(\n  context,\n  \"avec ces coefficients puisque ces arbres que nous allons cr\u00e9er\",\n  \"avec ces coefficients puisque les arbres que nous allons cr\u00e9er\"\n);\n\n// 10. \"leurs efficacit\u00e9s\" -> \"leur efficacit\u00e9\"\nawait replaceOnce(\n  context,\n  \"dans le but d\u2019augmenter leurs efficacit\u00e9s\",\n  \"dans le but d\u2019augmenter leur efficacit\u00e9\"\n);\n\n// 11. add comma after \"groupe\"\nawait replaceOnce(\n  context,\n  \"Malo \u00e9tant le meilleur programmeur du groupe il se chargera\",\n  \"Malo \u00e9tant le meilleur programmeur du groupe, il se chargera\"\n);\n\n// 12. \"g\u00e9rer\" -> \"g\u00e9r\u00e9e\", \"d'alaise\" -> \"\u00e0 l'aise\"\nawait replaceOnce(\n  context,\n  \"Cette partie \u00e9tant g\u00e9rer par quelqu\u2019un d\u2019alaise \",\n  \"Cette partie \u00e9tant g\u00e9r\u00e9e par quelqu\u2019un \u00e0 l\u2019aise \"\n);\n\n// 13. \"d\u00e9cid\u00e9s\" -> \"d\u00e9cid\u00e9\"\nawait replaceOnce(\n  context,\n  \"nous avons d\u00e9cid\u00e9s de placer deux personnes pour ce travail\",\n  \"nous avons d\u00e9cid\u00e9 de placer deux personnes pour ce travail\"\n);\n\n// 14. \"pr\u00e9vus\" -> \"pr\u00e9vu\"\nawait replaceOnce(\n  context,\n  \"Nous avons pr\u00e9vus au moins un rendez-vous\",\n  \"Nous avons pr\u00e9vu au moins un rendez-vous\"\n);\n\n// 15. add comma after \"deux\"\nawait replaceOnce(\n  context,\n  \"les s\u00e9ances de projet un mercredi sur deux ce qui nous semble\",\n  \"les s\u00e9ances de projet un mercredi sur deux, ce qui nous semble\"\n);\n\n// 16. \"affich\u00e9s, cr\u00e9er\" -> \"affich\u00e9s. Cr\u00e9er\"\nawait replaceOnce(\n  context,\n  \"les contrats et les arbres peuvent maintenant \u00eatre affich\u00e9s, cr\u00e9er\",\n  \"les contrats et les arbres peuvent maintenant \u00eatre affich\u00e9s. Cr\u00e9er\"\n);\n", "ps1": "# Applies the textual corrections described in the commit diff to the\n# synopsis document using the Word COM object model (Find/Replace on\n# $d.Content, anchored on unambiguous surrounding context).\n\n$d = $word.ActiveDocument\n\nfunction Replace-OnceInDoc($findText, $replaceText) {\n    $rng = $d.Content\n    $find = $rng.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $result = $find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 1)\n    if (-not $result) {\n        throw \"Text not found: $findText\"\n    }\n}\n\n# 1. \"ressource\" -> \"ressources\"\nReplace-OnceInDoc \"un nombre limit\u00e9 de ressource qu\u2019il devra g\u00e9rer\" \"un nombre limit\u00e9 de ressources qu\u2019il devra g\u00e9rer\"\n\n# 2. add commas around \"qui servira pour bon nombre de choses\"\nReplace-OnceInDoc \"g\u00e9rer son budget qui servira pour bon nombre de choses \u00e0 savoir marchander\" \"g\u00e9rer son budget, qui servira pour bon nombre de choses, \u00e0 savoir marchander\"\n\n# 3. add comma after \"croissante\"\nReplace-OnceInDoc \"le logement d\u2019une population croissante ou encore pour la gestion\" \"le logement d\u2019une population croissante, ou encore pour la gestion\"\n\n# 4. remove duplicate \"l\u00e0\"\nReplace-OnceInDoc \"mais \u00e9galement l\u00e0 pour vous racheter\" \"mais \u00e9galement pour vous racheter\"\n\n# 5. \"tel que\" -> \"telles que\" + commas\nReplace-OnceInDoc \"au cours du jeu tel que l\u2019arbre de comp\u00e9tences octroyant des bonus non n\u00e9gligeables ou encore un syst\u00e8me\" \"au cours du jeu, telles que l\u2019arbre de comp\u00e9tences octroyant des bonus non n\u00e9gligeables, ou encore un syst\u00e8me\"\n\n# 6. \"car dans\" -> \"car, dans\"\nReplace-OnceInDoc \" car dans le cas contraire cela veut dire\" \" car, dans le cas contraire cela veut dire\"\n\n# 7. \"ne doivent pas\" -> \"de\"\nReplace-OnceInDoc \"Emp\u00eacher toutes les IA de la carte ne doivent pas faire la m\u00eame chose\" \"Emp\u00eacher toutes les IA de la carte de faire la m\u00eame chose\"\n\n# 8. \"etc\" -> \"etc.\"\nReplace-OnceInDoc \"ressources etc + arbres de comp\u00e9tences\" \"ressources etc. + arbres de comp\u00e9tences\"\n\n# 9. \"ces arbres\" -> \"les arbres\"\nReplace-OnceInDoc \"avec ces coefficients puisque ces arbres que nous allons cr\u00e9er\" \"avec ces coefficients puisque les arbres que nous allons cr\u00e9er\"\n\n# 10. \"leurs efficacit\u00e9s\" -> \"leur efficacit\u00e9\"\nReplace-OnceInDoc \"dans le but d\u2019augmenter leurs efficacit\u00e9s\" \"dans le but d\u2019augmenter leur efficacit\u00e9\"\n\n# 11. add comma after \"groupe\"\nReplace-OnceInDoc \"Malo \u00e9tant le meilleur programmeur du groupe il se chargera\" \"Malo \u00e9tant le meilleur programmeur du groupe, il se chargera\"\n\n# 12. \"g\u00e9rer\" -> \"g\u00e9r\u00e9e\", \"d\u2019alaise\" -> \"\u00e0 l\u2019aise\"\nReplace-OnceInDoc \"Cette partie \u00e9tant g\u00e9rer par quelqu\u2019un d\u2019alaise \" \"Cette partie \u00e9tant g\u00e9r\u00e9e par quelqu\u2019un \u00e0 l\u2019aise \"\n\n# 13. \"d\u00e9cid\u00e9s\" -> \"d\u00e9cid\u00e9\"\nReplace-OnceInDoc \"nous avons d\u00e9cid\u00e9s de placer deux personnes pour ce travail\" \"nous avons d\u00e9cid\u00e9 de placer deux personnes pour ce travail\"\n\n# 14. \"pr\u00e9vus\" -> \"pr\u00e9vu\"\nReplace-OnceInDoc \"Nous avons pr\u00e9vus au moins un rendez-vous\" \"Nous avons pr\u00e9vu au moins un rendez-vous\"\n\n# 15. add comma after \"deux\"\nReplace-OnceInDoc \"les s\u00e9ances de projet un mercredi sur deux ce qui nous semble\" \"les s\u00e9ances de projet un mercredi sur deux, ce qui nous semble\"\n\n# 16. \"affich\u00e9s, cr\u00e9er\" -> \"affich\u00e9s. Cr\u00e9er\"\nReplace-OnceInDoc \"les contrats et les arbres peuvent maintenant \u00eatre affich\u00e9s, cr\u00e9er\" \"les contrats et les arbres peuvent maintenant \u00eatre affich\u00e9s. Cr\u00e9er\"\n"}
